$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Countries/provincias data refresh ---------------------------------
# The sheet is sorted by "Casos totales" (col B) descending. Bangladesh's
# case count jumped by 969 (row currently at A38), which leapfrogs it past
# Poland, Israel, Ukraine, Austria and Japan (rows 33-37). Israel also picks
# up a small update (+20 new cases) in the same refresh.
#
# Update the changed countries' figures first, then re-sort the A33:H38
# block by total cases (desc) so Bangladesh's row floats up to the top and
# Poland/Israel/Ukraine/Austria/Japan shift down one row each, exactly as
# in the refreshed dashboard.

# Bangladesh (currently row 38): new totals
$ws.Range("B38").Value = 16660
$ws.Range("C38").Value = 969
$ws.Range("D38").Value = 3147
$ws.Range("E38").Value = 13263
$ws.Range("F38").Value = 1
$ws.Range("G38").Value = 11
$ws.Range("H38").Value = 250

# Israel (currently row 34): new totals
$ws.Range("B34").Value = 16526
$ws.Range("C34").Value = 20
$ws.Range("D34").Value = 11956
$ws.Range("E34").Value = 4312
$ws.Range("F34").Value = 67
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 258

# Re-sort so the updated Bangladesh row takes its new, higher rank.
$rng = $ws.Range("A33:H38")
$rng.Sort($ws.Range("B33"), 2)

# Dinamarca (row 46): new totals (recoveries/critical untouched)
$ws.Range("B46").Value = 10591
$ws.Range("C46").Value = 78
$ws.Range("E46").Value = 1730

# Footer timestamp refresh
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 11:05"
